# Commit: "add sort for authors and objects"
#  - B3: update the commit date range in the report header text
#  - The per-(type) blocks of data rows are re-sorted by the "object"
#    column (C), using a stable sort so rows that share the same object
#    keep their original relative order (i.e. authors are not reordered
#    within an object - "sort for ... objects", ties broken by original
#    row order which is effectively the author insertion order).
#
# NOTE: Sort-Object in this host is NOT a stable sort once the pipeline
# has more than a couple of items (verified empirically), so we roll our
# own stable insertion sort over plain PowerShell arrays instead of
# relying on it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Header text: date range changed from 2016-01-01 to 2010-01-01
$ws.Range("B3").Value = "Коммиты с 2010-01-01 по 2024-03-01"

function StableSortByObjectColumn($list) {
    # $list: array of row-arrays, each row-array is 1-based-friendly
    # (we always index with [1] for column "object" = column C, which is
    # offset 1 within the B:H block we read).
    $n = $list.Count
    $indices = New-Object 'object[]' $n
    for ($i = 0; $i -lt $n; $i++) { $indices[$i] = $i }

    # Simple stable insertion sort on the index array.
    for ($i = 1; $i -lt $n; $i++) {
        $j = $i
        while ($j -gt 0) {
            $a = $list[$indices[$j - 1]][1]
            $b = $list[$indices[$j]][1]
            if ($a -gt $b) {
                $tmp = $indices[$j - 1]
                $indices[$j - 1] = $indices[$j]
                $indices[$j] = $tmp
                $j--
            } else {
                break
            }
        }
    }

    $result = @()
    foreach ($idx in $indices) { $result += , $list[$idx] }
    return $result
}

function SortBlock($rangeAddr) {
    $rng = $ws.Range($rangeAddr)
    $arr = $rng.Value2
    $nrows = $arr.GetLength(0)
    $ncols = $arr.GetLength(1)

    $rows = @()
    for ($i = 1; $i -le $nrows; $i++) {
        $rowVals = @()
        for ($j = 1; $j -le $ncols; $j++) {
            $rowVals += $arr[$i, $j]
        }
        $rows += , $rowVals
    }

    $sorted = StableSortByObjectColumn $rows

    # NB: New-Object 'object[,]' is 0-based (unlike the COM Value2 array,
    # which is 1-based on read) - build it 0-based here.
    $outArr = New-Object 'object[,]' $nrows, $ncols
    for ($i = 0; $i -lt $nrows; $i++) {
        $r = $sorted[$i]
        for ($j = 0; $j -lt $ncols; $j++) {
            $outArr[$i, $j] = $r[$j]
        }
    }
    $rng.Value2 = $outArr
}

# 2) Re-sort the CommonModules rows (7-34) by object (column C)
SortBlock "B7:H34"

# 3) Re-sort the Catalogs rows (45-57) by object (column C)
SortBlock "B45:H57"
